$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Cells.Item(2, 2).Value = -0.053409393884569618
$ws.Cells.Item(2, 3).Value = 0.010936365897328697
$ws.Cells.Item(2, 4).Value = -0.074844305747975653
$ws.Cells.Item(2, 5).Value = -0.031974482021163583

$ws.Cells.Item(3, 2).Value = -0.041197781007235118
$ws.Cells.Item(3, 3).Value = 0.0082237504040701698
$ws.Cells.Item(3, 4).Value = -0.057316057697845438
$ws.Cells.Item(3, 5).Value = -0.025079504316624795

$ws.Cells.Item(4, 2).Value = -0.058353280975296344
$ws.Cells.Item(4, 3).Value = 0.010978920009276026
$ws.Cells.Item(4, 4).Value = -0.079871603514824224
$ws.Cells.Item(4, 5).Value = -0.036834958435768464

$ws.Cells.Item(5, 2).Value = -0.052446826787112946
$ws.Cells.Item(5, 3).Value = 0.0076706208043054933
$ws.Cells.Item(5, 4).Value = -0.067480991189071193
$ws.Cells.Item(5, 5).Value = -0.037412662385154699

$ws.Cells.Item(6, 2).Value = -0.050648497471115293
$ws.Cells.Item(6, 3).Value = 0.025752187017136972
$ws.Cells.Item(6, 4).Value = -0.10112224713685772
$ws.Cells.Item(6, 5).Value = -0.00017474780537286155

$ws.Cells.Item(7, 2).Value = 0.03746401207669995
$ws.Cells.Item(7, 3).Value = 0.030113397796505244
$ws.Cells.Item(7, 4).Value = -0.021557757149188263
$ws.Cells.Item(7, 5).Value = 0.096485781302588169

$ws.Cells.Item(8, 2).Value = -0.023950716954273775
$ws.Cells.Item(8, 3).Value = 0.0099203725303353921
$ws.Cells.Item(8, 4).Value = -0.043394315754477886
$ws.Cells.Item(8, 5).Value = -0.0045071181540696607

$ws.Cells.Item(9, 2).Value = -0.08011349332521682
$ws.Cells.Item(9, 3).Value = 0.0078654733145022207
$ws.Cells.Item(9, 4).Value = -0.095529558861813441
$ws.Cells.Item(9, 5).Value = -0.064697427788620199

$ws.Cells.Item(10, 2).Value = -0.026268242600263537
$ws.Cells.Item(10, 3).Value = 0.0095261684227990698
$ws.Cells.Item(10, 4).Value = -0.044939219755982272
$ws.Cells.Item(10, 5).Value = -0.0075972654445448026

$ws.Cells.Item(11, 2).Value = -0.091024337669478331
$ws.Cells.Item(11, 3).Value = 0.0074142282705403646
$ws.Cells.Item(11, 4).Value = -0.10555598114094084
$ws.Cells.Item(11, 5).Value = -0.076492694198015826

$ws.Cells.Item(12, 2).Value = -0.048632434415056117
$ws.Cells.Item(12, 3).Value = 0.018031871169262161
$ws.Cells.Item(12, 4).Value = -0.08397452597398275
$ws.Cells.Item(12, 5).Value = -0.013290342856129483

$ws.Cells.Item(13, 2).Value = 0.010617891994827608
$ws.Cells.Item(13, 3).Value = 0.022808629971091079
$ws.Cells.Item(13, 4).Value = -0.034086651265602688
$ws.Cells.Item(13, 5).Value = 0.055322435255257901

$ws.Cells.Item(14, 2).Value = -0.022392717653295573
$ws.Cells.Item(14, 3).Value = 0.0087675705571767758
$ws.Cells.Item(14, 4).Value = -0.039576863091880618
$ws.Cells.Item(14, 5).Value = -0.0052085722147105322

$ws.Cells.Item(15, 2).Value = -0.070864650549751676
$ws.Cells.Item(15, 3).Value = 0.0081235451890629348
$ws.Cells.Item(15, 4).Value = -0.08678652835883352
$ws.Cells.Item(15, 5).Value = -0.054942772740669832

$ws.Cells.Item(16, 2).Value = -0.025038654029724099
$ws.Cells.Item(16, 3).Value = 0.01011747190063377
$ws.Cells.Item(16, 4).Value = -0.044868566576541535
$ws.Cells.Item(16, 5).Value = -0.0052087414829066654

$ws.Cells.Item(17, 2).Value = -0.071257780999121709
$ws.Cells.Item(17, 3).Value = 0.0082486387399135792
$ws.Cells.Item(17, 4).Value = -0.087424841537248596
$ws.Cells.Item(17, 5).Value = -0.055090720460994821

$ws.Cells.Item(18, 2).Value = -0.042853518449864424
$ws.Cells.Item(18, 3).Value = 0.015997487959408203
$ws.Cells.Item(18, 4).Value = -0.07420826133079228
$ws.Cells.Item(18, 5).Value = -0.01149877556893656

$ws.Cells.Item(19, 2).Value = -0.033860754281792839
$ws.Cells.Item(19, 3).Value = 0.018278719709737042
$ws.Cells.Item(19, 4).Value = -0.069686747208369759
$ws.Cells.Item(19, 5).Value = 0.0019652386447840883

$ws.Cells.Item(20, 2).Value = -0.015678345003433111
$ws.Cells.Item(20, 3).Value = 0.0079783100179282272
$ws.Cells.Item(20, 4).Value = -0.031315566147889579
$ws.Cells.Item(20, 5).Value = -0.000041123858976640365

$ws.Cells.Item(21, 2).Value = -0.07393945256542428
$ws.Cells.Item(21, 3).Value = 0.0077502965155812117
$ws.Cells.Item(21, 4).Value = -0.089129775415031171
$ws.Cells.Item(21, 5).Value = -0.058749129715817382

$ws.Cells.Item(22, 2).Value = -0.028574021678970718
$ws.Cells.Item(22, 3).Value = 0.0089586527422834891
$ws.Cells.Item(22, 4).Value = -0.046132686744861304
$ws.Cells.Item(22, 5).Value = -0.011015356613080136

$ws.Cells.Item(23, 2).Value = -0.070166369964854911
$ws.Cells.Item(23, 3).Value = 0.0082931240325458538
$ws.Cells.Item(23, 4).Value = -0.086420620212910818
$ws.Cells.Item(23, 5).Value = -0.053912119716799003

$ws.Cells.Item(24, 2).Value = -0.03078374369269184
$ws.Cells.Item(24, 3).Value = 0.014818303200214676
$ws.Cells.Item(24, 4).Value = -0.059827309029530484
$ws.Cells.Item(24, 5).Value = -0.0017401783558531961

$ws.Cells.Item(25, 2).Value = -0.052110823145036421
$ws.Cells.Item(25, 3).Value = 0.019076548059742276
$ws.Cells.Item(25, 4).Value = -0.089500546643436929
$ws.Cells.Item(25, 5).Value = -0.014721099646635913

$ws.Cells.Item(26, 2).Value = -0.01161340344436318
$ws.Cells.Item(26, 3).Value = 0.007955625820075956
$ws.Cells.Item(26, 4).Value = -0.02720616431872442
$ws.Cells.Item(26, 5).Value = 0.0039793574299980608

$ws.Cells.Item(27, 2).Value = -0.078420509604172517
$ws.Cells.Item(27, 3).Value = 0.0074577348740955387
$ws.Cells.Item(27, 4).Value = -0.093037421387680994
$ws.Cells.Item(27, 5).Value = -0.06380359782066404

$ws.Cells.Item(28, 2).Value = -0.02077853139611166
$ws.Cells.Item(28, 3).Value = 0.0083935017921517707
$ws.Cells.Item(28, 4).Value = -0.037229519166040517
$ws.Cells.Item(28, 5).Value = -0.0043275436261828032

$ws.Cells.Item(29, 2).Value = -0.081977128516457026
$ws.Cells.Item(29, 3).Value = 0.0072084501501968861
$ws.Cells.Item(29, 4).Value = -0.096105453642443078
$ws.Cells.Item(29, 5).Value = -0.067848803390470974

$ws.Cells.Item(30, 2).Value = -0.029859602011330557
$ws.Cells.Item(30, 3).Value = 0.013949939663068655
$ws.Cells.Item(30, 4).Value = -0.057201192919216812
$ws.Cells.Item(30, 5).Value = -0.0025180111034443023

$ws.Cells.Item(31, 2).Value = -0.084045460015076096
$ws.Cells.Item(31, 3).Value = 0.019187553904729576
$ws.Cells.Item(31, 4).Value = -0.12165275316170558
$ws.Cells.Item(31, 5).Value = -0.046438166868446609

$ws.Cells.Item(32, 2).Value = -0.0059324932926342853
$ws.Cells.Item(32, 3).Value = 0.0068220497077040622
$ws.Cells.Item(32, 4).Value = -0.019303482850338224
$ws.Cells.Item(32, 5).Value = 0.007438496265069653

$ws.Cells.Item(33, 2).Value = -0.066179524663568884
$ws.Cells.Item(33, 3).Value = 0.0072198034928526774
$ws.Cells.Item(33, 4).Value = -0.080330098870203245
$ws.Cells.Item(33, 5).Value = -0.052028950456934515

$ws.Cells.Item(34, 2).Value = -0.0053613654916186059
$ws.Cells.Item(34, 3).Value = 0.0088356851972075899
$ws.Cells.Item(34, 4).Value = -0.022679018208880956
$ws.Cells.Item(34, 5).Value = 0.011956287225643743

$ws.Cells.Item(35, 2).Value = -0.072348590643591962
$ws.Cells.Item(35, 3).Value = 0.0079813405404029567
$ws.Cells.Item(35, 4).Value = -0.087991755505174074
$ws.Cells.Item(35, 5).Value = -0.056705425782009856

$ws.Cells.Item(36, 2).Value = -0.048215388669748525
$ws.Cells.Item(36, 3).Value = 0.01343621422937853
$ws.Cells.Item(36, 4).Value = -0.074550088437874681
$ws.Cells.Item(36, 5).Value = -0.021880688901622373

$ws.Cells.Item(37, 2).Value = -0.091815541213802962
$ws.Cells.Item(37, 3).Value = 0.017250635239779013
$ws.Cells.Item(37, 4).Value = -0.12562650532361766
$ws.Cells.Item(37, 5).Value = -0.058004577103988271

$ws.Cells.Item(38, 2).Value = -0.014340747297011494
$ws.Cells.Item(38, 3).Value = 0.0071080862575655362
$ws.Cells.Item(38, 4).Value = -0.028272358938280784
$ws.Cells.Item(38, 5).Value = -0.00040913565574220677

$ws.Cells.Item(39, 2).Value = -0.052588511744695046
$ws.Cells.Item(39, 3).Value = 0.0080106199255186028
$ws.Cells.Item(39, 4).Value = -0.068289059801080756
$ws.Cells.Item(39, 5).Value = -0.036887963688309336

$ws.Cells.Item(40, 2).Value = 0.00336118063851725
$ws.Cells.Item(40, 3).Value = 0.0095446722326958996
$ws.Cells.Item(40, 4).Value = -0.015346063376713502
$ws.Cells.Item(40, 5).Value = 0.022068424653748003

$ws.Cells.Item(41, 2).Value = -0.047777721789058959
$ws.Cells.Item(41, 3).Value = 0.0079707371342824171
$ws.Cells.Item(41, 4).Value = -0.063400104323512146
$ws.Cells.Item(41, 5).Value = -0.032155339254605765

$ws.Cells.Item(42, 2).Value = -0.05646317500588266
$ws.Cells.Item(42, 3).Value = 0.012720735913217679
$ws.Cells.Item(42, 4).Value = -0.08139555219079915
$ws.Cells.Item(42, 5).Value = -0.031530797820966169

$ws.Cells.Item(43, 2).Value = -0.10446731836023264
$ws.Cells.Item(43, 3).Value = 0.020290147169746153
$ws.Cells.Item(43, 4).Value = -0.14423567634841844
$ws.Cells.Item(43, 5).Value = -0.064698960372046838

$ws.Cells.Item(44, 2).Value = -0.010327060341784782
$ws.Cells.Item(44, 3).Value = 0.0086082311652528696
$ws.Cells.Item(44, 4).Value = -0.027198905894436791
$ws.Cells.Item(44, 5).Value = 0.0065447852108672291

$ws.Cells.Item(45, 2).Value = -0.052156658587050771
$ws.Cells.Item(45, 3).Value = 0.0078350534908732342
$ws.Cells.Item(45, 4).Value = -0.067513102283241228
$ws.Cells.Item(45, 5).Value = -0.036800214890860314

$ws.Cells.Item(46, 2).Value = 0.0039297776527214322
$ws.Cells.Item(46, 3).Value = 0.010457811890258763
$ws.Cells.Item(46, 4).Value = -0.016567190092937289
$ws.Cells.Item(46, 5).Value = 0.024426745398380151

$ws.Cells.Item(47, 2).Value = -0.039904707690273801
$ws.Cells.Item(47, 3).Value = 0.0084735362297295967
$ws.Cells.Item(47, 4).Value = -0.056512559908988513
$ws.Cells.Item(47, 5).Value = -0.023296855471559089

$ws.Cells.Item(48, 2).Value = -0.048255270848543488
$ws.Cells.Item(48, 3).Value = 0.014161312926831399
$ws.Cells.Item(48, 4).Value = -0.076011148946643611
$ws.Cells.Item(48, 5).Value = -0.020499392750443369

$ws.Cells.Item(49, 2).Value = -0.12046451497794451
$ws.Cells.Item(49, 3).Value = 0.021343723783655094
$ws.Cells.Item(49, 4).Value = -0.16229786596984719
$ws.Cells.Item(49, 5).Value = -0.078631163986041835

$ws.Cells.Item(50, 2).Value = 0.0010442271643383485
$ws.Cells.Item(50, 3).Value = 0.0090871050615831186
$ws.Cells.Item(50, 4).Value = -0.016766195229824782
$ws.Cells.Item(50, 5).Value = 0.018854649558501482

$ws.Cells.Item(51, 2).Value = -0.04117040178923792
$ws.Cells.Item(51, 3).Value = 0.010150435114437709
$ws.Cells.Item(51, 4).Value = -0.061064916294884856
$ws.Cells.Item(51, 5).Value = -0.021275887283590984

$ws.Cells.Item(52, 2).Value = 0.0010867476377361843
$ws.Cells.Item(52, 3).Value = 0.010105098767109646
$ws.Cells.Item(52, 4).Value = -0.018718913973855215
$ws.Cells.Item(52, 5).Value = 0.020892409249327587

$ws.Cells.Item(53, 2).Value = -0.027689940664026299
$ws.Cells.Item(53, 3).Value = 0.012292349666584455
$ws.Cells.Item(53, 4).Value = -0.051782541574510219
$ws.Cells.Item(53, 5).Value = -0.0035973397535423828

$ws.Cells.Item(54, 2).Value = -0.018751182657079361
$ws.Cells.Item(54, 3).Value = 0.017947173004440455
$ws.Cells.Item(54, 4).Value = -0.053927267578763351
$ws.Cells.Item(54, 5).Value = 0.016424902264604628

$ws.Cells.Item(55, 2).Value = -0.086288820328780727
$ws.Cells.Item(55, 3).Value = 0.02534047128157824
$ws.Cells.Item(55, 4).Value = -0.13595573132180497
$ws.Cells.Item(55, 5).Value = -0.036621909335756479
